# Scheduled data refresh: updates current market-board pricing snapshots
# (currentAveragePrice* / LevePrice* / LeveProfit* columns) for the affected
# Leve rows across the crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 5000
$ws.Range("K62").Value = 5000
$ws.Range("M62").Value = -4376
# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 5000
$ws.Range("K65").Value = 25000
$ws.Range("M65").Value = -21880
# Row 107: Another Man's Ink
$ws.Range("H107").Value = 423.81818
$ws.Range("I107").Value = 446.2
$ws.Range("K107").Value = 446.2
$ws.Range("M107").Value = 1473.8
# Row 116: Growing Up
$ws.Range("H116").Value = 9634.333000000001
$ws.Range("I116").Value = 9456.75
$ws.Range("K116").Value = 9456.75
$ws.Range("M116").Value = -6014.75
# Row 138: All-night Crafting
$ws.Range("H138").Value = 3247.5781
$ws.Range("I138").Value = 3031
$ws.Range("J138").Value = 3278.5178
$ws.Range("K138").Value = 9093
$ws.Range("L138").Value = 9835.553400000001
$ws.Range("M138").Value = -3953
$ws.Range("N138").Value = -20115.5534

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 1148.2858
$ws.Range("J2").Value = 1370.6
$ws.Range("L2").Value = 1370.6
$ws.Range("N2").Value = -1596.6
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 9088.052
$ws.Range("I32").Value = 6272.3228
$ws.Range("J32").Value = 19999
$ws.Range("K32").Value = 6272.3228
$ws.Range("L32").Value = 19999
$ws.Range("M32").Value = -5985.3228
$ws.Range("N32").Value = -20573
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 2870.6
$ws.Range("I61").Value = 2838.875
$ws.Range("K61").Value = 2838.875
$ws.Range("M61").Value = -2626.875
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 799.25
$ws.Range("I74").Value = 799.25
$ws.Range("K74").Value = 799.25
$ws.Range("M74").Value = 74.75
# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 799.25
$ws.Range("I77").Value = 799.25
$ws.Range("K77").Value = 3996.25
$ws.Range("M77").Value = 371.75
# Row 102: Smells of Rich Tama-hagane
$ws.Range("H102").Value = 2312.7
$ws.Range("I102").Value = 2312.7
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2312.7
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -690.6999999999998
$ws.Range("N102").ClearContents()
# Row 116: No Scope
$ws.Range("H116").Value = 1148.2858
$ws.Range("J116").Value = 1370.6
$ws.Range("L116").Value = 1370.6
$ws.Range("N116").Value = -5958.6
# Row 136: Metal with Mettle
$ws.Range("H136").Value = 2870.6
$ws.Range("I136").Value = 2838.875
$ws.Range("K136").Value = 8516.625
$ws.Range("M136").Value = -5966.625

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 1148.2858
$ws.Range("J3").Value = 1370.6
$ws.Range("L3").Value = 1370.6
$ws.Range("N3").Value = -1598.6
# Row 94: High Steal
$ws.Range("H94").Value = 1421.25
$ws.Range("I94").Value = 784
$ws.Range("K94").Value = 784
$ws.Range("M94").Value = -333

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 1234.8125
$ws.Range("I31").Value = 1159.1538
$ws.Range("J31").Value = 1562.6666
$ws.Range("K31").Value = 1159.1538
$ws.Range("L31").Value = 1562.6666
$ws.Range("M31").Value = -864.1538
$ws.Range("N31").Value = -2152.6666
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 1234.8125
$ws.Range("I34").Value = 1159.1538
$ws.Range("J34").Value = 1562.6666
$ws.Range("K34").Value = 1159.1538
$ws.Range("L34").Value = 1562.6666
$ws.Range("M34").Value = -957.1538
$ws.Range("N34").Value = -1966.6666
# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 3695.4285
$ws.Range("I132").Value = 3344.25
$ws.Range("J132").Value = 3835.9
$ws.Range("K132").Value = 10032.75
$ws.Range("L132").Value = 11507.7
$ws.Range("M132").Value = -7502.75
$ws.Range("N132").Value = -16567.7
# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 3325.3333
$ws.Range("I134").Value = 3325.3333
$ws.Range("K134").Value = 9975.999899999999
$ws.Range("M134").Value = -7440.999899999999

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 6599.2
$ws.Range("I80").Value = 3992.5
$ws.Range("K80").Value = 3992.5
$ws.Range("M80").Value = -2994.5
# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 6599.2
$ws.Range("I83").Value = 3992.5
$ws.Range("K83").Value = 19962.5
$ws.Range("M83").Value = -14970.5
# Row 92: Play It by Ear
$ws.Range("H92").Value = 38333.332
$ws.Range("J92").Value = 38333.332
$ws.Range("L92").Value = 38333.332
$ws.Range("N92").Value = -42077.332
# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 965.2857
$ws.Range("I122").Value = 1271.6
$ws.Range("K122").Value = 3814.8
$ws.Range("M122").Value = -1364.8
# Row 126: Gold Rush Order
$ws.Range("H126").Value = 1490
$ws.Range("I126").Value = 1490
$ws.Range("K126").Value = 4470
$ws.Range("M126").Value = -2000
# Row 132: On Board for Lar
$ws.Range("H132").Value = 3189.1333
$ws.Range("I132").Value = 2557.3333
$ws.Range("J132").Value = 3610.3333
$ws.Range("K132").Value = 7671.999899999999
$ws.Range("L132").Value = 10830.9999
$ws.Range("M132").Value = -5141.999899999999
$ws.Range("N132").Value = -15890.9999
# Row 135: Fan of the Foreign
$ws.Range("H135").Value = 75000
$ws.Range("J135").Value = 75000
$ws.Range("L135").Value = 75000
$ws.Range("N135").Value = -85140

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 4699.5
$ws.Range("J22").Value = 5110.5557
$ws.Range("L22").Value = 5110.5557
$ws.Range("N22").Value = -5700.5557
# Row 27: Fire and Hide
$ws.Range("H27").Value = 4699.5
$ws.Range("J27").Value = 5110.5557
$ws.Range("L27").Value = 5110.5557
$ws.Range("N27").Value = -5324.5557
# Row 46: Supply Side Logic
$ws.Range("H46").Value = 2985
$ws.Range("I46").Value = 2977.5
$ws.Range("K46").Value = 2977.5
$ws.Range("M46").Value = -2789.5
# Row 82: Trainin' the Neck
$ws.Range("H82").Value = 976
$ws.Range("I82").Value = 976
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 976
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -615
$ws.Range("N82").ClearContents()
# Row 85: Training Is Only Skintight (L)
$ws.Range("H85").Value = 976
$ws.Range("I85").Value = 976
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 976
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 272
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 104: Brimming with Confidence
$ws.Range("H104").Value = 26650
$ws.Range("J104").Value = 26650
$ws.Range("L104").Value = 26650
$ws.Range("N104").Value = -33638
# Row 107: Flax Wax
$ws.Range("H107").Value = 598.8889
$ws.Range("I107").Value = 547.2857
$ws.Range("J107").Value = 779.5
$ws.Range("K107").Value = 1641.8571
$ws.Range("L107").Value = 2338.5
$ws.Range("M107").Value = 278.1428999999998
$ws.Range("N107").Value = -6178.5
# Row 122: Heavy Armoire
$ws.Range("H122").Value = 3102.3845
$ws.Range("I122").Value = 3110.9167
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 9332.750100000001
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -6882.750100000001
$ws.Range("N122").Value = -13900
# Row 125: Color Coated
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 4356.05
$ws.Range("I132").Value = 3998
$ws.Range("J132").Value = 4714.1
$ws.Range("K132").Value = 11994
$ws.Range("L132").Value = 14142.3
$ws.Range("M132").Value = -9464
$ws.Range("N132").Value = -19202.3
